# Updates cryptos list data (price + volume columns, plus a couple of
# row reorderings) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the cell to stay a plain text value (matches the inlineStr cells
    # already present in the sheet) instead of letting Excel autoconvert
    # numeric-looking strings (e.g. "1.00", "0.999") into numbers, and reset
    # the style afterwards so no stray quote-prefix / text-format style sticks.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "62.943.48"
Set-TextCell "E2" "  -5.42%  "
Set-TextCell "D3" "3.102.02"
Set-TextCell "E3" "  -6.14%  "
Set-TextCell "E4" "  -0.04%  "
Set-TextCell "D5" "555.65"
Set-TextCell "E5" "  -5.37%  "
Set-TextCell "D6" "160.95"
Set-TextCell "E6" "  -11.12%  "
Set-TextCell "E7" "  +0.03%  "
Set-TextCell "D8" "0.582"
Set-TextCell "E8" "  -10.01%  "
Set-TextCell "D9" "3.108.65"
Set-TextCell "E9" "  -6.02%  "
Set-TextCell "D10" "6.72"
Set-TextCell "E10" "  -2.01%  "
Set-TextCell "D11" "0.114"
Set-TextCell "E11" "  -9.50%  "
Set-TextCell "D12" "0.373"
Set-TextCell "E12" "  -7.53%  "
Set-TextCell "D13" "3.654.26"
Set-TextCell "E13" "  -5.99%  "
Set-TextCell "E14" "  -2.03%  "
Set-TextCell "D15" "63.011.95"
Set-TextCell "E15" "  -5.30%  "
Set-TextCell "D16" "24.33"
Set-TextCell "E16" "  -8.73%  "
Set-TextCell "D17" "3.110.57"
Set-TextCell "E17" "  -4.28%  "
Set-TextCell "D18" "0.0000151"
Set-TextCell "E18" "  -7.86%  "
Set-TextCell "D19" "400.59"
Set-TextCell "E19" "  -6.11%  "
Set-TextCell "D20" "12.34"
Set-TextCell "E20" "  -5.66%  "
Set-TextCell "D21" "5.11"
Set-TextCell "E21" "  -6.64%  "
Set-TextCell "D22" "7.00"
Set-TextCell "E22" "  -4.68%  "
Set-TextCell "D23" "0.999"
Set-TextCell "E23" "  -0.08%  "
Set-TextCell "E24" "  -0.02%  "
Set-TextCell "D25" "67.75"
Set-TextCell "E25" "  -5.47%  "
Set-TextCell "D26" "0.199"
Set-TextCell "E26" "  -3.79%  "
Set-TextCell "D27" "0.480"
Set-TextCell "E27" "  -6.90%  "
Set-TextCell "D28" "0.0000100"
Set-TextCell "E28" "  -12.72%  "
Set-TextCell "D29" "8.60"
Set-TextCell "E29" "  -5.53%  "
Set-TextCell "D30" "0.999"
Set-TextCell "E30" "  -0.09%  "
Set-TextCell "E31" "  -0.05%  "
Set-TextCell "D32" "1.76"
Set-TextCell "E32" "  -8.05%  "
Set-TextCell "D33" "20.82"
Set-TextCell "E33" "  -6.94%  "
Set-TextCell "D34" "4.77"
Set-TextCell "E34" "  -7.68%  "
Set-TextCell "D35" "6.13"
Set-TextCell "E35" "  -7.00%  "
Set-TextCell "D36" "152.47"
Set-TextCell "E36" "  -4.73%  "
Set-TextCell "D37" "1.08"
Set-TextCell "E37" "  -8.53%  "
Set-TextCell "D38" "1.31"
Set-TextCell "E38" "  -8.86%  "
Set-TextCell "D39" "2.697.00"
Set-TextCell "E39" "  -6.21%  "
Set-TextCell "D40" "1.62"
Set-TextCell "E40" "  -9.77%  "
Set-TextCell "D41" "23.32"
Set-TextCell "E41" "  -11.69%  "
Set-TextCell "D42" "4.01"
Set-TextCell "E42" "  -7.42%  "
Set-TextCell "D43" "38.06"
Set-TextCell "E43" "  -4.31%  "
Set-TextCell "E44" "  -8.07%  "
Set-TextCell "D45" "0.0597"
Set-TextCell "E45" "  -9.28%  "
Set-TextCell "B46" "RenderToken"
Set-TextCell "C46" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D46" "5.19"
Set-TextCell "E46" "  -12.73%  "
Set-TextCell "B47" "VeChain"
Set-TextCell "C47" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D47" "0.0252"
Set-TextCell "E47" "  -7.24%  "
Set-TextCell "D48" "281.68"
Set-TextCell "E48" "  -9.26%  "
Set-TextCell "B49" "FirstDigitalUSD"
Set-TextCell "C49" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D49" "1.00"
Set-TextCell "E49" "  -0.03%  "
Set-TextCell "B50" "InjectiveProtocol"
Set-TextCell "C50" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D50" "20.50"
Set-TextCell "E50" "  -11.04%  "
Set-TextCell "D51" "0.0959"
Set-TextCell "E51" "  -7.04%  "
